$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 8.850858
$ws.Cells.Item(2, 8).Value = 26.552574
$ws.Cells.Item(2, 9).Value = 0.05442939716240135
$ws.Cells.Item(2, 10).Value = 0.05442939716240137
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.065175
$ws.Cells.Item(2, 14).Value = 0.195525
$ws.Cells.Item(2, 15).Value = 0.009404016458916581
$ws.Cells.Item(2, 16).Value = 0.009404016458916581
$ws.Cells.Item(2, 17).Value = 0.57685467015
$ws.Cells.Item(2, 18).Value = 5.19169203135
$ws.Cells.Item(2, 19).Value = 0.0005118549467641298
$ws.Cells.Item(2, 20).Value = 0.0005118549467641299

$ws.Cells.Item(3, 7).Value = 8.850858
$ws.Cells.Item(3, 8).Value = 26.552574
$ws.Cells.Item(3, 9).Value = 0.05442939716240135
$ws.Cells.Item(3, 10).Value = 0.05442939716240137
$ws.Cells.Item(3, 13).Value = 6.718514333333332
$ws.Cells.Item(3, 14).Value = 20.155543
$ws.Cells.Item(3, 15).Value = 0.969405744075698
$ws.Cells.Item(3, 16).Value = 0.969405744075698
$ws.Cells.Item(3, 17).Value = 59.464616335298
$ws.Cells.Item(3, 18).Value = 535.181547017682
$ws.Cells.Item(3, 19).Value = 0.05276417025580937
$ws.Cells.Item(3, 20).Value = 0.05276417025580939

$ws.Cells.Item(4, 7).Value = 8.850858
$ws.Cells.Item(4, 8).Value = 26.552574
$ws.Cells.Item(4, 9).Value = 0.05442939716240135
$ws.Cells.Item(4, 10).Value = 0.05442939716240137
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.14686
$ws.Cells.Item(4, 14).Value = 0.44058
$ws.Cells.Item(4, 15).Value = 0.02119023946538534
$ws.Cells.Item(4, 16).Value = 0.02119023946538533
$ws.Cells.Item(4, 17).Value = 1.29983700588
$ws.Cells.Item(4, 18).Value = 11.69853305292
$ws.Cells.Item(4, 19).Value = 0.00115337195982785
$ws.Cells.Item(4, 20).Value = 0.00115337195982785

$ws.Cells.Item(5, 9).Value = 0.2412735821509021
$ws.Cells.Item(5, 10).Value = 0.2412735821509022
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.065175
$ws.Cells.Item(5, 14).Value = 0.195525
$ws.Cells.Item(5, 15).Value = 0.009404016458916581
$ws.Cells.Item(5, 16).Value = 0.009404016458916581
$ws.Cells.Item(5, 17).Value = 2.55707025805
$ws.Cells.Item(5, 18).Value = 23.01363232245
$ws.Cells.Item(5, 19).Value = 0.002268940737648845
$ws.Cells.Item(5, 20).Value = 0.002268940737648846

$ws.Cells.Item(6, 9).Value = 0.2412735821509021
$ws.Cells.Item(6, 10).Value = 0.2412735821509022
$ws.Cells.Item(6, 13).Value = 6.718514333333332
$ws.Cells.Item(6, 14).Value = 20.155543
$ws.Cells.Item(6, 15).Value = 0.969405744075698
$ws.Cells.Item(6, 16).Value = 0.969405744075698
$ws.Cells.Item(6, 17).Value = 263.5936046037482
$ws.Cells.Item(6, 18).Value = 2372.342441433734
$ws.Cells.Item(6, 19).Value = 0.2338919964308043
$ws.Cells.Item(6, 20).Value = 0.2338919964308044

$ws.Cells.Item(7, 9).Value = 0.2412735821509021
$ws.Cells.Item(7, 10).Value = 0.2412735821509022
$ws.Cells.Item(7, 11).Value = 1
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.14686
$ws.Cells.Item(7, 14).Value = 0.44058
$ws.Cells.Item(7, 15).Value = 0.02119023946538534
$ws.Cells.Item(7, 16).Value = 0.02119023946538533
$ws.Cells.Item(7, 17).Value = 5.761892414226668
$ws.Cells.Item(7, 18).Value = 51.85703172804001
$ws.Cells.Item(7, 19).Value = 0.005112644982448937
$ws.Cells.Item(7, 20).Value = 0.005112644982448939

$ws.Cells.Item(8, 7).Value = 50.430027
$ws.Cells.Item(8, 8).Value = 151.290081
$ws.Cells.Item(8, 9).Value = 0.3101254102702387
$ws.Cells.Item(8, 10).Value = 0.3101254102702387
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.065175
$ws.Cells.Item(8, 14).Value = 0.195525
$ws.Cells.Item(8, 15).Value = 0.009404016458916581
$ws.Cells.Item(8, 16).Value = 0.009404016458916581
$ws.Cells.Item(8, 17).Value = 3.286777009725
$ws.Cells.Item(8, 18).Value = 29.58099308752501
$ws.Cells.Item(8, 19).Value = 0.002916424462509582
$ws.Cells.Item(8, 20).Value = 0.002916424462509582

$ws.Cells.Item(9, 7).Value = 50.430027
$ws.Cells.Item(9, 8).Value = 151.290081
$ws.Cells.Item(9, 9).Value = 0.3101254102702387
$ws.Cells.Item(9, 10).Value = 0.3101254102702387
$ws.Cells.Item(9, 13).Value = 6.718514333333332
$ws.Cells.Item(9, 14).Value = 20.155543
$ws.Cells.Item(9, 15).Value = 0.969405744075698
$ws.Cells.Item(9, 16).Value = 0.969405744075698
$ws.Cells.Item(9, 17).Value = 338.814859229887
$ws.Cells.Item(9, 18).Value = 3049.333733068983
$ws.Cells.Item(9, 19).Value = 0.3006373540998018
$ws.Cells.Item(9, 20).Value = 0.3006373540998019

$ws.Cells.Item(10, 7).Value = 50.430027
$ws.Cells.Item(10, 8).Value = 151.290081
$ws.Cells.Item(10, 9).Value = 0.3101254102702387
$ws.Cells.Item(10, 10).Value = 0.3101254102702387
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.14686
$ws.Cells.Item(10, 14).Value = 0.44058
$ws.Cells.Item(10, 15).Value = 0.02119023946538534
$ws.Cells.Item(10, 16).Value = 0.02119023946538533
$ws.Cells.Item(10, 17).Value = 7.406153765220001
$ws.Cells.Item(10, 18).Value = 66.65538388698
$ws.Cells.Item(10, 19).Value = 0.00657163170792723
$ws.Cells.Item(10, 20).Value = 0.00657163170792723

$ws.Cells.Item(11, 7).Value = 3.269985333333333
$ws.Cells.Item(11, 8).Value = 9.809956
$ws.Cells.Item(11, 9).Value = 0.02010916121614733
$ws.Cells.Item(11, 10).Value = 0.02010916121614734
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.065175
$ws.Cells.Item(11, 14).Value = 0.195525
$ws.Cells.Item(11, 15).Value = 0.009404016458916581
$ws.Cells.Item(11, 16).Value = 0.009404016458916581
$ws.Cells.Item(11, 17).Value = 0.2131212941
$ws.Cells.Item(11, 18).Value = 1.9180916469
$ws.Cells.Item(11, 19).Value = 0.0001891068830516565
$ws.Cells.Item(11, 20).Value = 0.0001891068830516566

$ws.Cells.Item(12, 7).Value = 3.269985333333333
$ws.Cells.Item(12, 8).Value = 9.809956
$ws.Cells.Item(12, 9).Value = 0.02010916121614733
$ws.Cells.Item(12, 10).Value = 0.02010916121614734
$ws.Cells.Item(12, 13).Value = 6.718514333333332
$ws.Cells.Item(12, 14).Value = 20.155543
$ws.Cells.Item(12, 15).Value = 0.969405744075698
$ws.Cells.Item(12, 16).Value = 0.969405744075698
$ws.Cells.Item(12, 17).Value = 21.96944333178977
$ws.Cells.Item(12, 18).Value = 197.724989986108
$ws.Cells.Item(12, 19).Value = 0.01949393639147748
$ws.Cells.Item(12, 20).Value = 0.01949393639147748

$ws.Cells.Item(13, 7).Value = 3.269985333333333
$ws.Cells.Item(13, 8).Value = 9.809956
$ws.Cells.Item(13, 9).Value = 0.02010916121614733
$ws.Cells.Item(13, 10).Value = 0.02010916121614734
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.14686
$ws.Cells.Item(13, 14).Value = 0.44058
$ws.Cells.Item(13, 15).Value = 0.02119023946538534
$ws.Cells.Item(13, 16).Value = 0.02119023946538533
$ws.Cells.Item(13, 17).Value = 0.4802300460533334
$ws.Cells.Item(13, 18).Value = 4.32207041448
$ws.Cells.Item(13, 19).Value = 0.0004261179416182014
$ws.Cells.Item(13, 20).Value = 0.0004261179416182015

$ws.Cells.Item(14, 7).Value = 47.074941
$ws.Cells.Item(14, 8).Value = 141.224823
$ws.Cells.Item(14, 9).Value = 0.2894929124482182
$ws.Cells.Item(14, 10).Value = 0.2894929124482182
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.065175
$ws.Cells.Item(14, 14).Value = 0.195525
$ws.Cells.Item(14, 15).Value = 0.009404016458916581
$ws.Cells.Item(14, 16).Value = 0.009404016458916581
$ws.Cells.Item(14, 17).Value = 3.068109279675
$ws.Cells.Item(14, 18).Value = 27.612983517075
$ws.Cells.Item(14, 19).Value = 0.002722396113402741
$ws.Cells.Item(14, 20).Value = 0.002722396113402741

$ws.Cells.Item(15, 7).Value = 47.074941
$ws.Cells.Item(15, 8).Value = 141.224823
$ws.Cells.Item(15, 9).Value = 0.2894929124482182
$ws.Cells.Item(15, 10).Value = 0.2894929124482182
$ws.Cells.Item(15, 13).Value = 6.718514333333332
$ws.Cells.Item(15, 14).Value = 20.155543
$ws.Cells.Item(15, 15).Value = 0.969405744075698
$ws.Cells.Item(15, 16).Value = 0.969405744075698
$ws.Cells.Item(15, 17).Value = 316.273665849321
$ws.Cells.Item(15, 18).Value = 2846.462992643889
$ws.Cells.Item(15, 19).Value = 0.2806360921965058
$ws.Cells.Item(15, 20).Value = 0.2806360921965059

$ws.Cells.Item(16, 7).Value = 47.074941
$ws.Cells.Item(16, 8).Value = 141.224823
$ws.Cells.Item(16, 9).Value = 0.2894929124482182
$ws.Cells.Item(16, 10).Value = 0.2894929124482182
$ws.Cells.Item(16, 11).Value = 1
$ws.Cells.Item(16, 12).Value = 0.3333333333333333
$ws.Cells.Item(16, 13).Value = 0.14686
$ws.Cells.Item(16, 14).Value = 0.44058
$ws.Cells.Item(16, 15).Value = 0.02119023946538534
$ws.Cells.Item(16, 16).Value = 0.02119023946538533
$ws.Cells.Item(16, 17).Value = 6.913425835260001
$ws.Cells.Item(16, 18).Value = 62.22083251734001
$ws.Cells.Item(16, 19).Value = 0.006134424138309575
$ws.Cells.Item(16, 20).Value = 0.006134424138309575

$ws.Cells.Item(17, 7).Value = 13.751998
$ws.Cells.Item(17, 8).Value = 41.255994
$ws.Cells.Item(17, 9).Value = 0.08456953675209218
$ws.Cells.Item(17, 10).Value = 0.0845695367520922
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.065175
$ws.Cells.Item(17, 14).Value = 0.195525
$ws.Cells.Item(17, 15).Value = 0.009404016458916581
$ws.Cells.Item(17, 16).Value = 0.009404016458916581
$ws.Cells.Item(17, 17).Value = 0.89628646965
$ws.Cells.Item(17, 18).Value = 8.06657822685
$ws.Cells.Item(17, 19).Value = 0.0007952933155396255
$ws.Cells.Item(17, 20).Value = 0.0007952933155396258

$ws.Cells.Item(18, 7).Value = 13.751998
$ws.Cells.Item(18, 8).Value = 41.255994
$ws.Cells.Item(18, 9).Value = 0.08456953675209218
$ws.Cells.Item(18, 10).Value = 0.0845695367520922
$ws.Cells.Item(18, 13).Value = 6.718514333333332
$ws.Cells.Item(18, 14).Value = 20.155543
$ws.Cells.Item(18, 15).Value = 0.969405744075698
$ws.Cells.Item(18, 16).Value = 0.969405744075698
$ws.Cells.Item(18, 17).Value = 92.39299567497132
$ws.Cells.Item(18, 18).Value = 831.536961074742
$ws.Cells.Item(18, 19).Value = 0.08198219470129901
$ws.Cells.Item(18, 20).Value = 0.08198219470129903

$ws.Cells.Item(19, 7).Value = 13.751998
$ws.Cells.Item(19, 8).Value = 41.255994
$ws.Cells.Item(19, 9).Value = 0.08456953675209218
$ws.Cells.Item(19, 10).Value = 0.0845695367520922
$ws.Cells.Item(19, 11).Value = 1
$ws.Cells.Item(19, 12).Value = 0.3333333333333333
$ws.Cells.Item(19, 13).Value = 0.14686
$ws.Cells.Item(19, 14).Value = 0.44058
$ws.Cells.Item(19, 15).Value = 0.02119023946538534
$ws.Cells.Item(19, 16).Value = 0.02119023946538533
$ws.Cells.Item(19, 17).Value = 2.01961842628
$ws.Cells.Item(19, 18).Value = 18.17656583652
$ws.Cells.Item(19, 19).Value = 0.00179204873525354
$ws.Cells.Item(19, 20).Value = 0.00179204873525354
